$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 1038, shifting existing rows 1038-1090 down to 1041-1093
$ws.Rows.Item(1038).Resize(3).Insert()

# Populate the 3 newly inserted rows with the new Hass / Peru price records
# Row 1038
$ws.Cells.Item(1038, 1).Value = 5
$ws.Cells.Item(1038, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1038, 3).Value = "Maule"
$ws.Cells.Item(1038, 4).Value = 44746
$ws.Cells.Item(1038, 5).Value = 7
$ws.Cells.Item(1038, 6).Value = "Fruta"
$ws.Cells.Item(1038, 7).Value = 100106
$ws.Cells.Item(1038, 8).Value = "Oleaginosos"
$ws.Cells.Item(1038, 9).Value = 100106002
$ws.Cells.Item(1038, 10).Value = "Palta"
$ws.Cells.Item(1038, 11).Value = "Hass"
$ws.Cells.Item(1038, 12).Value = "3a nueva (o)"
$ws.Cells.Item(1038, 13).Value = 400
$ws.Cells.Item(1038, 14).Value = 10000
$ws.Cells.Item(1038, 15).Value = 10000
$ws.Cells.Item(1038, 16).Value = 10000
$ws.Cells.Item(1038, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(1038, 18).Value = "Perú"
$ws.Cells.Item(1038, 19).Value = 1000
$ws.Cells.Item(1038, 20).Value = 10

# Row 1039
$ws.Cells.Item(1039, 1).Value = 5
$ws.Cells.Item(1039, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1039, 3).Value = "Maule"
$ws.Cells.Item(1039, 4).Value = 44746
$ws.Cells.Item(1039, 5).Value = 7
$ws.Cells.Item(1039, 6).Value = "Fruta"
$ws.Cells.Item(1039, 7).Value = 100106
$ws.Cells.Item(1039, 8).Value = "Oleaginosos"
$ws.Cells.Item(1039, 9).Value = 100106002
$ws.Cells.Item(1039, 10).Value = "Palta"
$ws.Cells.Item(1039, 11).Value = "Hass"
$ws.Cells.Item(1039, 12).Value = "Segunda"
$ws.Cells.Item(1039, 13).Value = 850
$ws.Cells.Item(1039, 14).Value = 15000
$ws.Cells.Item(1039, 15).Value = 15000
$ws.Cells.Item(1039, 16).Value = 15000
$ws.Cells.Item(1039, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(1039, 18).Value = "Perú"
$ws.Cells.Item(1039, 19).Value = 1500
$ws.Cells.Item(1039, 20).Value = 10

# Row 1040
$ws.Cells.Item(1040, 1).Value = 5
$ws.Cells.Item(1040, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1040, 3).Value = "Maule"
$ws.Cells.Item(1040, 4).Value = 44746
$ws.Cells.Item(1040, 5).Value = 7
$ws.Cells.Item(1040, 6).Value = "Fruta"
$ws.Cells.Item(1040, 7).Value = 100106
$ws.Cells.Item(1040, 8).Value = "Oleaginosos"
$ws.Cells.Item(1040, 9).Value = 100106002
$ws.Cells.Item(1040, 10).Value = "Palta"
$ws.Cells.Item(1040, 11).Value = "Hass"
$ws.Cells.Item(1040, 12).Value = "Tercera"
$ws.Cells.Item(1040, 13).Value = 500
$ws.Cells.Item(1040, 14).Value = 13000
$ws.Cells.Item(1040, 15).Value = 13000
$ws.Cells.Item(1040, 16).Value = 13000
$ws.Cells.Item(1040, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(1040, 18).Value = "Perú"
$ws.Cells.Item(1040, 19).Value = 1300
$ws.Cells.Item(1040, 20).Value = 10
